$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("classFields")

$ws.Range("B2").Value = "`$VALUES"
$ws.Range("C2").Value = "private"
$ws.Range("D2").Value = "org.andante.enums.KafkaConsumerGroup[]"
$ws.Range("B3").Value = "ACTIVITY_PRODUCT_GROUP"
$ws.Range("C3").Value = "public"
$ws.Range("D3").Value = "org.andante.enums.KafkaConsumerGroup"
$ws.Range("B4").Value = "name"
$ws.Range("C4").Value = "private"
$ws.Range("D4").Value = "java.lang.String"
$ws.Range("B5").Value = "PRODUCT_ORDER_GROUP"
$ws.Range("C5").Value = "public"
$ws.Range("D5").Value = "org.andante.enums.KafkaConsumerGroup"
$ws.Range("B6").Value = "ACTIVITY_ORDER_GROUP"
$ws.Range("B7").Value = "eurekaContainer"
$ws.Range("D7").Value = "org.testcontainers.containers.GenericContainer"
$ws.Range("B8").Value = "postgresqlContainer"
$ws.Range("D8").Value = "org.testcontainers.containers.PostgreSQLContainer"
$ws.Range("B9").Value = "kafkaContainer"
$ws.Range("D9").Value = "org.testcontainers.containers.KafkaContainer"
$ws.Range("B12").Value = "NOT_FOUND"
$ws.Range("C12").Value = "public"
$ws.Range("D12").Value = "org.andante.enums.OperationStatus"
$ws.Range("B13").Value = "`$VALUES"
$ws.Range("C13").Value = "private"
$ws.Range("D13").Value = "org.andante.enums.OperationStatus[]"
$ws.Range("B16").Value = "PRODUCT_HEADPHONES_VARIANT_INTERNAL_TOPIC"
$ws.Range("B17").Value = "PRODUCT_SUBWOOFERS_INTERNAL_TOPIC"
$ws.Range("B18").Value = "PRODUCT_SPEAKERS_VARIANT_INTERNAL_TOPIC"
$ws.Range("B19").Value = "ORDER_ENTRY_INTERNAL_TOPIC"
$ws.Range("B20").Value = "PRODUCT_MICROPHONE_VARIANT_INTERNAL_TOPIC"
$ws.Range("B21").Value = "PRODUCT_PRODUCER_INTERNAL_TOPIC"
$ws.Range("B22").Value = "topicName"
$ws.Range("C22").Value = "private"
$ws.Range("D22").Value = "java.lang.String"
$ws.Range("B23").Value = "PRODUCT_GRAMOPHONE_VARIANT_INTERNAL_TOPIC"
$ws.Range("B24").Value = "PRODUCT_COMMENT_INTERNAL_TOPIC"
$ws.Range("B25").Value = "ORDER_INTERNAL_TOPIC"
$ws.Range("B26").Value = "PRODUCT_SPEAKERS_INTERNAL_TOPIC"
$ws.Range("B27").Value = "PRODUCT_MICROPHONE_INTERNAL_TOPIC"
$ws.Range("B28").Value = "PRODUCT_ORDER_INTERNAL_TOPIC"
$ws.Range("C28").Value = "public"
$ws.Range("D28").Value = "org.andante.enums.KafkaTopic"
$ws.Range("B30").Value = "PRODUCT_HEADPHONES_INTERNAL_TOPIC"
$ws.Range("B31").Value = "PRODUCT_AMPLIFIER_INTERNAL_TOPIC"
$ws.Range("C31").Value = "public"
$ws.Range("D31").Value = "org.andante.enums.KafkaTopic"
$ws.Range("B32").Value = "`$VALUES"
$ws.Range("C32").Value = "private"
$ws.Range("D32").Value = "org.andante.enums.KafkaTopic[]"
$ws.Range("B33").Value = "PRODUCT_SUBWOOFERS_VARIANT_INTERNAL_TOPIC"
$ws.Range("B34").Value = "serialVersionUID"
$ws.Range("D34").Value = "long"
$ws.Range("B35").Value = "`$assertionsDisabled"
$ws.Range("D35").Value = "boolean"
$ws.Range("B37").Value = "`$VALUES"
$ws.Range("D37").Value = "org.andante.rsql.operator.RSQLSearchOperator[]"
$ws.Range("B38").Value = "IN"
$ws.Range("B39").Value = "NOT_IN"
$ws.Range("B41").Value = "operator"
$ws.Range("C41").Value = "private"
$ws.Range("D41").Value = "cz.jirutka.rsql.parser.ast.ComparisonOperator"
$ws.Range("B42").Value = "GREATER_THAN_OR_EQUAL"
$ws.Range("B43").Value = "LESS_THAN"
$ws.Range("B44").Value = "LESS_THAN_OR_EQUAL"
$ws.Range("C44").Value = "public"
$ws.Range("D44").Value = "org.andante.rsql.operator.RSQLSearchOperator"
$ws.Range("B45").Value = "NOT_EQUAL"
$ws.Range("B46").Value = "GREATER_THAN"
$ws.Range("B49").Value = "`$VALUES"
$ws.Range("C49").Value = "private"
$ws.Range("D49").Value = "org.andante.enums.OperationType[]"
$ws.Range("B50").Value = "CREATE"
$ws.Range("B51").Value = "MODIFY"
$ws.Range("B52").Value = "DELETE"
$ws.Range("C52").Value = "public"
$ws.Range("D52").Value = "org.andante.enums.OperationType"
$ws.Range("B53").Value = "property"
$ws.Range("B54").Value = "operator"
$ws.Range("D54").Value = "cz.jirutka.rsql.parser.ast.ComparisonOperator"
$ws.Range("B55").Value = "RSQL_MAPPING_ERROR_MESSAGE"
$ws.Range("D55").Value = "java.lang.String"
$ws.Range("B56").Value = "DATABASE_WILDCARD"
$ws.Range("D56").Value = "java.lang.String"
$ws.Range("B57").Value = "RSQL_WILDCARD"
$ws.Range("B58").Value = "arguments"
$ws.Range("D58").Value = "java.util.List"
$ws.Range("B60").Value = "`$VALUES"
$ws.Range("C60").Value = "private"
$ws.Range("D60").Value = "org.andante.rsql.operator.RSQLReservedOperator[]"
$ws.Range("B61").Value = "LESS_THAN"
$ws.Range("C61").Value = "public"
$ws.Range("D61").Value = "org.andante.rsql.operator.RSQLReservedOperator"
$ws.Range("B62").Value = "APOSTROPHE"
$ws.Range("B63").Value = "TILDE"
$ws.Range("B64").Value = "QUOTE"
$ws.Range("B65").Value = "EXCLAMATION_MARK"
$ws.Range("B67").Value = "RIGHT_BRACKET"
$ws.Range("B68").Value = "reservedCharacter"
$ws.Range("C68").Value = "private"
$ws.Range("D68").Value = "java.lang.String"
$ws.Range("B69").Value = "LEFT_BRACKET"
$ws.Range("B70").Value = "EQUALS"
$ws.Range("C70").Value = "public"
$ws.Range("D70").Value = "org.andante.rsql.operator.RSQLReservedOperator"
$ws.Range("B71").Value = "GREATER_THAN"
$ws.Range("B72").Value = "SEMICOLON"
